# The post "「奮闘なくして進歩なし」" (row 740) was removed from the sheet.
# Deleting the entire row shifts every subsequent row up by one, which
# matches the renumbering seen throughout the rest of the sheet (old row
# 741 -> new row 740, ..., old row 844 -> new row 843) and updates the
# sheet's used-range dimension from A1:C844 to A1:C843 automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(740).Delete()
